# Weekly update: insert a new price-record row for "Feria Lagunitas de
# Puerto Montt - Brócoli" at row 436, pushing the existing rows 436-457
# down to 437-458 (dimension grows from A1:R457 to A1:R458).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 436 (existing row 436 and below shift down by one).
$ws.Rows.Item(436).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A436").Value = 4
$ws.Range("B436").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C436").Value = "Los Lagos"
$ws.Range("D436").Value = 44931
$ws.Range("E436").Value = 10
$ws.Range("F436").Value = 100112023
$ws.Range("G436").Value = "Brócoli"
$ws.Range("H436").Value = "Sin especificar"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 500
$ws.Range("K436").Value = 1500
$ws.Range("L436").Value = 1500
$ws.Range("M436").Value = 1500
$ws.Range("N436").Value = '$/unidad'
$ws.Range("O436").Value = "Región Metropolitana"
$ws.Range("P436").Value = 1500
$ws.Range("Q436").Value = 1
$ws.Range("R436").Value = "Hortaliza"
